# "Test Results" sheet -> "value type settings page added"
#
# The report now reflects a run from 28-12-2024 (was 27-12-2024), the 4th
# test case is renamed to "verifyValueTypeSettingMap" (replacing
# "verifyCustomerNavigationAfterLogin"), and the trailing duplicate rows
# for "verifyTipsCategoryValueListAddition" (rows 5-8) are removed, leaving
# a 4-row table (header + 3 results).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the execution date for every test row: 27-12-2024 -> 28-12-2024
$ws.Cells.Replace("27-12-2024", "28-12-2024")

# 2) Rename the 4th test case
$ws.Cells.Item(4, 1).Value = "verifyValueTypeSettingMap"

# 3) Drop the extra "verifyTipsCategoryValueListAddition" rows (5 through 8)
$ws.Rows("5:8").Delete()

# 4) Column A's best-fit width shrinks now that the longest method name changed
$ws.Columns.Item(1).AutoFit()

Write-Host "value type settings page added: sheet updated"
